{"js": "// Replace the arithmetic-answer text in every table cell according to the\n// old-answer -> new-answer mapping derived from the commit's diff.\n// (Table is a 20-row x 5-col grid of \"A+B=C\" / \"A-B=C\" strings; every\n// old value in the map is unique in the document, so a straight lookup\n// by current cell text is unambiguous and robust to any cell reordering.)\nconst answerMap = {\"14+12=26\": \"27+42=69\", \"54-2=52\": \"74-36=38\", \"62+13=75\": \"73-66=7\", \"35+12=47\": \"49-48=1\", \"39-34=5\": \"13+11=24\", \"10+59=69\": \"2+88=90\", \"33-16=17\": \"28+28=56\", \"42-20=22\": \"30+44=74\", \"57-11=46\": \"52-24=28\", \"88-70=18\": \"59-14=45\", \"51-1=50\": \"74-72=2\", \"29+41=70\": \"39+12=51\", \"90-8=82\": \"65+16=81\", \"9+48=57\": \"94-71=23\", \"22+46=68\": \"91-67=24\", \"67+28=95\": \"78-41=37\", \"48+45=93\": \"12+78=90\", \"6+2=8\": \"54+36=90\", \"59-30=29\": \"48-24=24\", \"51-17=34\": \"71-41=30\", \"14-11=3\": \"31+17=48\", \"34+39=73\": \"86+1=87\", \"14+6=20\": \"30+32=62\", \"82-65=17\": \"26+70=96\", \"11+28=39\": \"75+19=94\", \"93-58=35\": \"72-8=64\", \"6+56=62\": \"58-52=6\", \"73+7=80\": \"87-11=76\", \"16+11=27\": \"55-11=44\", \"38-23=15\": \"65+34=99\", \"14+29=43\": \"10+47=57\", \"13+1=14\": \"81-31=50\", \"72+9=81\": \"69+10=79\", \"77-2=75\": \"41+23=64\", \"32+5=37\": \"42+48=90\", \"90-63=27\": \"42+13=55\", \"78-70=8\": \"23+75=98\", \"15+9=24\": \"23+67=90\", \"9+11=20\": \"41+45=86\", \"40+22=62\": \"96-60=36\", \"67+1=68\": \"79+2=81\", \"54+38=92\": \"19+62=81\", \"47-23=24\": \"38+42=80\", \"4+43=47\": \"22+40=62\", \"59-17=42\": \"94-15=79\", \"52+36=88\": \"37+10=47\", \"8+12=20\": \"3+29=32\", \"64-35=29\": \"89-22=67\", \"32-20=12\": \"56+23=79\", \"47+7=54\": \"33+43=76\", \"77-7=70\": \"56+13=69\", \"22+26=48\": \"44+36=80\", \"31+30=61\": \"72-52=20\", \"98-92=6\": \"29+42=71\", \"52-2=50\": \"98-87=11\", \"1+7=8\": \"83-41=42\", \"78-73=5\": \"46+53=99\", \"69-17=52\": \"69-15=54\", \"5+21=26\": \"4+5=9\", \"77+0=77\": \"42-40=2\", \"61+15=76\": \"86-76=10\", \"62-7=55\": \"77+21=98\", \"48+25=73\": \"34+18=52\", \"64-2=62\": \"65-55=10\", \"19+77=96\": \"54-54=0\", \"21+7=28\": \"92-89=3\", \"23-1=22\": \"8+1=9\", \"65+32=97\": \"67-16=51\", \"42+20=62\": \"6+11=17\", \"7+31=38\": \"30-24=6\", \"78-68=10\": \"20+56=76\", \"73-15=58\": \"39+15=54\", \"88-17=71\": \"33+32=65\", \"59-28=31\": \"59+19=78\", \"42+27=69\": \"74-71=3\", \"80-12=68\": \"96-17=79\", \"84+10=94\": \"15-6=9\", \"65-26=39\": \"89-56=33\", \"50-2=48\": \"2+92=94\", \"11+1=12\": \"2+81=83\", \"96-20=76\": \"39-26=13\", \"76-38=38\": \"34+25=59\", \"32+13=45\": \"98-32=66\", \"18+75=93\": \"18+26=44\", \"63+3=66\": \"33+63=96\", \"72+1=73\": \"54+30=84\", \"98-66=32\": \"66-6=60\", \"52-3=49\": \"19+43=62\", \"8+81=89\": \"10-7=3\", \"95-73=22\": \"95-74=21\", \"54-26=28\": \"10+46=56\", \"78+18=96\": \"51-46=5\", \"73-56=17\": \"20+38=58\", \"77-15=62\": \"19+65=84\", \"21+34=55\": \"92-23=69\", \"55+3=58\": \"7+15=22\", \"8+76=84\": \"89-42=47\", \"62-0=62\": \"6+5=11\", \"86-22=64\": \"8+56=64\", \"91-14=77\": \"18+7=25\"};\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Load every row's cells collection.\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\n// Load every cell's current value.\nconst allCells = [];\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    cell.load(\"value\");\n    allCells.push(cell);\n  }\n}\nawait context.sync();\n\n// Now assign the mapped replacement text onto each cell that needs it.\nfor (const cell of allCells) {\n  const current = cell.value;\n  if (Object.prototype.hasOwnProperty.call(answerMap, current)) {\n    cell.value = answerMap[current];\n  }\n}\nawait context.sync();\n", "ps1": "# Replace the arithmetic-answer text in every table cell according to the\n# old-answer -> new-answer mapping derived from the commit's diff.\n# (Table is a 20-row x 5-col grid of \"A+B=C\" / \"A-B=C\" strings; every\n# old value in the map is unique in the document, so a straight lookup\n# by current cell text is unambiguous and robust to any cell reordering.)\n$answerMap = @{\n  \"14+12=26\" = \"27+42=69\"\n  \"54-2=52\" = \"74-36=38\"\n  \"62+13=75\" = \"73-66=7\"\n  \"35+12=47\" = \"49-48=1\"\n  \"39-34=5\" = \"13+11=24\"\n  \"10+59=69\" = \"2+88=90\"\n  \"33-16=17\" = \"28+28=56\"\n  \"42-20=22\" = \"30+44=74\"\n  \"57-11=46\" = \"52-24=28\"\n  \"88-70=18\" = \"59-14=45\"\n  \"51-1=50\" = \"74-72=2\"\n  \"29+41=70\" = \"39+12=51\"\n  \"90-8=82\" = \"65+16=81\"\n  \"9+48=57\" = \"94-71=23\"\n  \"22+46=68\" = \"91-67=24\"\n  \"67+28=95\" = \"78-41=37\"\n  \"48+45=93\" = \"12+78=90\"\n  \"6+2=8\" = \"54+36=90\"\n  \"59-30=29\" = \"48-24=24\"\n  \"51-17=34\" = \"71-41=30\"\n  \"14-11=3\" = \"31+17=48\"\n  \"34+39=73\" = \"86+1=87\"\n  \"14+6=20\" = \"30+32=62\"\n  \"82-65=17\" = \"26+70=96\"\n  \"11+28=39\" = \"75+19=94\"\n  \"93-58=35\" = \"72-8=64\"\n  \"6+56=62\" = \"58-52=6\"\n  \"73+7=80\" = \"87-11=76\"\n  \"16+11=27\" = \"55-11=44\"\n  \"38-23=15\" = \"65+34=99\"\n  \"14+29=43\" = \"10+47=57\"\n  \"13+1=14\" = \"81-31=50\"\n  \"72+9=81\" = \"69+10=79\"\n  \"77-2=75\" = \"41+23=64\"\n  \"32+5=37\" = \"42+48=90\"\n  \"90-63=27\" = \"42+13=55\"\n  \"78-70=8\" = \"23+75=98\"\n  \"15+9=24\" = \"23+67=90\"\n  \"9+11=20\" = \"41+45=86\"\n  \"40+22=62\" = \"96-60=36\"\n  \"67+1=68\" = \"79+2=81\"\n  \"54+38=92\" = \"19+62=81\"\n  \"47-23=24\" = \"38+42=80\"\n  \"4+43=47\" = \"22+40=62\"\n  \"59-17=42\" = \"94-15=79\"\n  \"52+36=88\" = \"37+10=47\"\n  \"8+12=20\" = \"3+29=32\"\n  \"64-35=29\" = \"89-22=67\"\n  \"32-20=12\" = \"56+23=79\"\n  \"47+7=54\" = \"33+43=76\"\n  \"77-7=70\" = \"56+13=69\"\n  \"22+26=48\" = \"44+36=80\"\n  \"31+30=61\" = \"72-52=20\"\n  \"98-92=6\" = \"29+42=71\"\n  \"52-2=50\" = \"98-87=11\"\n  \"1+7=8\" = \"83-41=42\"\n  \"78-73=5\" = \"46+53=99\"\n  \"69-17=52\" = \"69-15=54\"\n  \"5+21=26\" = \"4+5=9\"\n  \"77+0=77\" = \"42-40=2\"\n  \"61+15=76\" = \"86-76=10\"\n  \"62-7=55\" = \"77+21=98\"\n  \"48+25=73\" = \"34+18=52\"\n  \"64-2=62\" = \"65-55=10\"\n  \"19+77=96\" = \"54-54=0\"\n  \"21+7=28\" = \"92-89=3\"\n  \"23-1=22\" = \"8+1=9\"\n  \"65+32=97\" = \"67-16=51\"\n  \"42+20=62\" = \"6+11=17\"\n  \"7+31=38\" = \"30-24=6\"\n  \"78-68=10\" = \"20+56=76\"\n  \"73-15=58\" = \"39+15=54\"\n  \"88-17=71\" = \"33+32=65\"\n  \"59-28=31\" = \"59+19=78\"\n  \"42+27=69\" = \"74-71=3\"\n  \"80-12=68\" = \"96-17=79\"\n  \"84+10=94\" = \"15-6=9\"\n  \"65-26=39\" = \"89-56=33\"\n  \"50-2=48\" = \"2+92=94\"\n  \"11+1=12\" = \"2+81=83\"\n  \"96-20=76\" = \"39-26=13\"\n  \"76-38=38\" = \"34+25=59\"\n  \"32+13=45\" = \"98-32=66\"\n  \"18+75=93\" = \"18+26=44\"\n  \"63+3=66\" = \"33+63=96\"\n  \"72+1=73\" = \"54+30=84\"\n  \"98-66=32\" = \"66-6=60\"\n  \"52-3=49\" = \"19+43=62\"\n  \"8+81=89\" = \"10-7=3\"\n  \"95-73=22\" = \"95-74=21\"\n  \"54-26=28\" = \"10+46=56\"\n  \"78+18=96\" = \"51-46=5\"\n  \"73-56=17\" = \"20+38=58\"\n  \"77-15=62\" = \"19+65=84\"\n  \"21+34=55\" = \"92-23=69\"\n  \"55+3=58\" = \"7+15=22\"\n  \"8+76=84\" = \"89-42=47\"\n  \"62-0=62\" = \"6+5=11\"\n  \"86-22=64\" = \"8+56=64\"\n  \"91-14=77\" = \"18+7=25\"\n}\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\nfor ($r = 1; $r -le $rowCount; $r++) {\n  for ($c = 1; $c -le $colCount; $c++) {\n    $cell = $t.Cell($r, $c)\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($answerMap.ContainsKey($current)) {\n      $cell.Range.Text = $answerMap[$current]\n    }\n  }\n}\n"}
